$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every cell as text (inline strings), including
# number-looking Price values such as '159.00' or '0.564'. Assigning a
# plain .Value lets Excel auto-detect those as numbers (dropping the
# trailing zero / changing the stored type), so for values that parse as
# a number we briefly force the Text format, assign, then restore the
# original General format.
$ws.Range('D2').Value = '61.652.30'
$ws.Range('E2').Value = '  -4.41%  '
$ws.Range('D3').Value = '2.966.97'
$ws.Range('E3').Value = '  -6.40%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '541.57'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -5.32%  '
$ws.Range('E6').Value = '  -7.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -3.12%  '
$ws.Range('D9').Value = '2.975.61'
$ws.Range('E9').Value = '  -6.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.113'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -4.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.12'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -7.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.367'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -4.49%  '
$ws.Range('D13').Value = '3.485.84'
$ws.Range('E13').Value = '  -6.40%  '
$ws.Range('E14').Value = '  -3.29%  '
$ws.Range('D15').Value = '61.687.46'
$ws.Range('E15').Value = '  -4.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.69'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -6.45%  '
$ws.Range('D17').Value = '2.972.99'
$ws.Range('E17').Value = '  -6.38%  '
$ws.Range('E18').Value = '  -5.40%  '
$ws.Range('E19').Value = '  -2.29%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '381.51'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -6.24%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.99'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -5.94%  '
$ws.Range('E22').Value = '  -6.82%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.15'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -5.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.471'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -3.36%  '
$ws.Range('D26').Value = '3.093.28'
$ws.Range('E26').Value = '  -6.33%  '
$ws.Range('E27').Value = '  -4.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').Value = '0.0₃0932'
$ws.Range('E29').Value = '  -8.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.32'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -5.76%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  -5.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.42'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '159.00'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.65'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.94'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -6.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.07'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('E38').Value = '  -5.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -8.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.92'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -4.40%  '
$ws.Range('D41').Value = '2.418.56'
$ws.Range('E41').Value = '  -10.20%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.28'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -7.42%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '37.08'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('E44').Value = '  -4.61%  '
$ws.Range('E45').Value = '  -4.31%  '
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0245'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -4.90%  '
$ws.Range('E48').Value = '  -9.79%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0952'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -3.30%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '267.88'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -8.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.73'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -8.09%  '
